# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet to make room for an additional "Variable Installments" style column,
# shifting the existing "Late" / "Outstanding (heading)" / "Outstanding"
# columns one place to the right, then select the sheet / cell that was
# being worked on when the change was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Remember the width of the column immediately to the left (M) so the
# freshly inserted column N can be given a matching width.
$leftWidth = $ws.Columns("M").ColumnWidth

# Insert a new blank column at N - this shifts the old N/O/P columns
# (and their data/styles) one column to the right, becoming O/P/Q.
$ws.Columns("N").Insert()

# Match the new column's width to its left neighbour (M).
$ws.Columns("N").ColumnWidth = $leftWidth

# Make "Repayment schedule" the active sheet and select cell K16, matching
# where the author was working when they saved the file.
$ws.Activate() | Out-Null
$ws.Range("K16").Select() | Out-Null
